$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")

$ws.Range("C2").Value = 0.07066352588744917
$ws.Range("D2").Value = 0.3901250515216378

$ws.Range("C3").Value = 0.1493128008527864
$ws.Range("D3").Value = 0.1342879180428377

$ws.Range("C4").Value = 0.0308465286395686
$ws.Range("D4").Value = 0.3770030367852095
